$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.511.66'
$ws.Range('E2').Value = '  +0.25%  '

$ws.Range('D3').Value = '2.290.15'
$ws.Range('E3').Value = '  +0.33%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').Value = "'504.05"
$ws.Range('E5').Value = '  +2.16%  '

$ws.Range('D6').Value = "'130.65"
$ws.Range('E6').Value = '  +3.02%  '

$ws.Range('D7').Value = "'0.997"
$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('D8').Value = "'0.531"
$ws.Range('E8').Value = '  +0.58%  '

$ws.Range('E9').Value = '  +1.82%  '

$ws.Range('E10').Value = '  +0.58%  '

$ws.Range('D11').Value = "'0.339"
$ws.Range('E11').Value = '  +5.30%  '

$ws.Range('D12').Value = "'4.76"
$ws.Range('E12').Value = '  +3.18%  '

$ws.Range('D13').Value = '2.702.17'
$ws.Range('E13').Value = '  +0.57%  '

$ws.Range('E14').Value = '  +6.60%  '

$ws.Range('D15').Value = '54.480.12'
$ws.Range('E15').Value = '  +0.35%  '

$ws.Range('E16').Value = '  +1.08%  '

$ws.Range('D17').Value = '2.303.15'
$ws.Range('E17').Value = '  +1.14%  '

$ws.Range('D18').Value = "'10.30"
$ws.Range('E18').Value = '  +3.38%  '

$ws.Range('D19').Value = "'4.18"
$ws.Range('E19').Value = '  +3.07%  '

$ws.Range('D20').Value = "'305.05"
$ws.Range('E20').Value = '  +0.77%  '

$ws.Range('E21').Value = '  -0.67%  '

$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('D23').Value = "'62.00"
$ws.Range('E23').Value = '  -2.50%  '

$ws.Range('D24').Value = "'0.996"
$ws.Range('E24').Value = '  -0.59%  '

$ws.Range('E25').Value = '  +1.88%  '

$ws.Range('E26').Value = '  +3.70%  '

$ws.Range('D27').Value = "'171.66"
$ws.Range('E27').Value = '  +1.85%  '

$ws.Range('E28').Value = '  +2.50%  '

$ws.Range('D29').Value = '0.0₃0696'
$ws.Range('E29').Value = '  +1.83%  '

$ws.Range('E30').Value = '  +1.61%  '

$ws.Range('E31').Value = '  +1.06%  '

$ws.Range('E32').Value = '  +0.15%  '

$ws.Range('E33').Value = '  +1.53%  '

$ws.Range('D34').Value = "'0.977"
$ws.Range('E34').Value = '  +12.07%  '

$ws.Range('D35').Value = "'0.996"
$ws.Range('E35').Value = '  -0.35%  '

$ws.Range('E36').Value = '  +0.82%  '

$ws.Range('D37').Value = "'3.75"
$ws.Range('E37').Value = '  +3.63%  '

$ws.Range('E38').Value = '  +0.70%  '

$ws.Range('E39').Value = '  +1.54%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = "'3.40"
$ws.Range('E40').Value = '  +1.90%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = "'4.89"
$ws.Range('E41').Value = '  +1.12%  '

$ws.Range('D42').Value = "'126.72"
$ws.Range('E42').Value = '  -0.51%  '

$ws.Range('E43').Value = '  +3.83%  '

$ws.Range('D44').Value = "'0.0900"
$ws.Range('E44').Value = '  +1.25%  '

$ws.Range('E45').Value = '  +1.36%  '

$ws.Range('D46').Value = "'243.02"
$ws.Range('E46').Value = '  +1.83%  '

$ws.Range('D47').Value = "'0.375"
$ws.Range('E47').Value = '  +0.65%  '

$ws.Range('D48').Value = "'0.0207"
$ws.Range('E48').Value = '  +2.05%  '

$ws.Range('E49').Value = '  +0.79%  '

$ws.Range('D50').Value = "'16.50"
$ws.Range('E50').Value = '  +1.47%  '

$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = "'1.54"
$ws.Range('E51').Value = '  +2.42%  '
